$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2225729989246953
$ws.Range("B3").Value = 0.2185938719739137
$ws.Range("B4").Value = 0.1580244899580968
$ws.Range("B5").Value = 0.1993691068915476
$ws.Range("B6").Value = 0.1919427087949038
$ws.Range("B7").Value = 0.2199786617192732
$ws.Range("B8").Value = 0.1762347010727719
$ws.Range("B9").Value = 0.1547639014287874
$ws.Range("B10").Value = 0.2515899659146517
$ws.Range("B11").Value = 0.2115133208638662
$ws.Range("B12").Value = 0.07994217227245783
$ws.Range("B13").Value = 0.172461673725689
$ws.Range("B14").Value = 0.1703472908511281
$ws.Range("B15").Value = 0.2065557988038792
$ws.Range("B16").Value = 0.1871266035253376
$ws.Range("B17").Value = 0.264368031372667
$ws.Range("B18").Value = 0.2034974490027841
$ws.Range("B19").Value = 0.1916350931379784
$ws.Range("B20").Value = 0.1780546825999253
$ws.Range("B21").Value = 0.196022225772513
$ws.Range("B22").Value = 0.1214537175050256
$ws.Range("B23").Value = 0.1533154264735849
$ws.Range("B24").Value = 0.1332916530968697
$ws.Range("B25").Value = 0.2028324594284665
$ws.Range("B26").Value = 0.1875431948010683
$ws.Range("B27").Value = 0.2577515256473847
$ws.Range("B28").Value = 0.3090118166696694
$ws.Range("B29").Value = 0.09961865905658623
$ws.Range("B30").Value = 0.1811816623748037
$ws.Range("B31").Value = 0.2826938213572558
$ws.Range("B32").Value = 0.1647730418862527
$ws.Range("B33").Value = 0.1209942430837056
$ws.Range("B34").Value = 0.1900194850035571
$ws.Range("B35").Value = 0.1703648828900786
$ws.Range("B36").Value = 0.1281530041310535
$ws.Range("B37").Value = 0.2341392414443393
$ws.Range("B38").Value = 0.2196548221740449
$ws.Range("B39").Value = 0.1834853024030706
$ws.Range("B40").Value = 0.2387031307088887
$ws.Range("B41").Value = 0.1613517033335533
$ws.Range("B42").Value = 0.1763442277687054
$ws.Range("B43").Value = 0.2828034277999462
$ws.Range("B44").Value = 0.2309564765118704
$ws.Range("B45").Value = 0.3246354056666368
$ws.Range("B46").Value = 0.1225666453223031
$ws.Range("B47").Value = 0.3235611708349139
$ws.Range("B48").Value = 0.1668956897754365
$ws.Range("B49").Value = 0.1851920853771044
$ws.Range("B50").Value = 0.1619949994789993
$ws.Range("B51").Value = 0.1813831568987349
$ws.Range("B52").Value = 0.2042228323230325
$ws.Range("B53").Value = 0.2098196311646337
$ws.Range("B54").Value = 0.1949395146105441
$ws.Range("B55").Value = 0.2430204226917999
$ws.Range("B56").Value = 0.1263052408444525
$ws.Range("B57").Value = 0.1121952593176505
$ws.Range("B58").Value = 0.1911460885545064
$ws.Range("B59").Value = 0.2065230419227556
$ws.Range("B60").Value = 0.330770000388631
$ws.Range("B61").Value = 0.2041524140415935
$ws.Range("B62").Value = 0.08326611516521817
$ws.Range("B63").Value = 0.1643457523645145
$ws.Range("B64").Value = 0.1968364552068901
$ws.Range("B65").Value = 0.1139138937921707
$ws.Range("B66").Value = 0.1482483440962012
$ws.Range("B67").Value = 0.1736791764473961
$ws.Range("B68").Value = 0.2591984599955779
$ws.Range("B69").Value = 0.1717094003668587
$ws.Range("B70").Value = 0.2117724920535764
$ws.Range("B71").Value = 0.1883695874210866
$ws.Range("B72").Value = 0.2566390044911464
$ws.Range("B73").Value = 0.2219694921722728
$ws.Range("B74").Value = 0.1456063870285668
$ws.Range("B75").Value = 0.1631545997309555
$ws.Range("B76").Value = 0.1174998185196909
$ws.Range("B77").Value = 0.1484457714358207
$ws.Range("B78").Value = 0.155349601895375
